$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# "Ready for handoff" rows (4-7) reflects the new handoff generation time.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-27 14:31:48"
$wsOverview.Range("G5").Value = "2016-08-27 14:31:48"
$wsOverview.Range("G6").Value = "2016-08-27 14:31:48"
$wsOverview.Range("G7").Value = "2016-08-27 14:31:48"

# zh-cn sheet: rows 4-7 (722dcc8b, c55049f1, c5a0839d, e884454c) just had a
# handoff generated -> Priority goes from "low" to "ht", and the handoff
# datetime is refreshed.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H4").Value = "2016-08-27 14:31:44"
$wsZhCn.Range("H5").Value = "2016-08-27 14:31:44"
$wsZhCn.Range("H6").Value = "2016-08-27 14:31:44"
$wsZhCn.Range("H7").Value = "2016-08-27 14:31:44"

# de-de sheet: same rows, same Priority change, and the handoff datetime
# is refreshed too.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("H4").Value = "2016-08-27 14:31:48"
$wsDeDe.Range("H5").Value = "2016-08-27 14:31:48"
$wsDeDe.Range("H6").Value = "2016-08-27 14:31:48"
$wsDeDe.Range("H7").Value = "2016-08-27 14:31:48"
